$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# Sheet 1: 展览 (Exhibition)
$ws1.Range("F3").Value = 857
$ws1.Range("F4").Value = 441
$ws1.Range("F5").Value = 889
$ws1.Range("F7").Value = 7580
$ws1.Range("F11").Value = 5495
$ws1.Range("F14").Value = 7573
$ws1.Range("F15").Value = 8925
$ws1.Range("F16").Value = 1142
$ws1.Range("F17").Value = 893
$ws1.Range("F18").Value = 4429
$ws1.Range("F19").Value = 667
$ws1.Range("F20").Value = 220
$ws1.Range("F21").Value = 83
$ws1.Range("B23").Value = "'2024-07-27"
$ws1.Range("C23").Value = "杭州·ACG CLUB动漫游戏嘉年华"
$ws1.Range("D23").Value = "中心路1号 白蓝地文创街区"
$ws1.Range("E23").Value = "2024.07.27 10:00-07.27 17:00"
$ws1.Range("F23").Value = 1192
$ws1.Range("G23").Value = 68.88
$ws1.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=86265"
$ws1.Range("I23").Value = "//i1.hdslb.com/bfs/openplatform/202405/XBRfeQwu1716533419093.jpeg"
$ws1.Range("C24").Value = "杭州·咒术回战ONLY2.0杀青宴"
$ws1.Range("D24").Value = "丰庆路492号建冠龙禾商务中心A幢 杭州华礼宴国际礼宴中心(龙禾商务中心店)"
$ws1.Range("E24").Value = "2024.07.27 10:00-07.27 19:00"
$ws1.Range("F24").Value = 110
$ws1.Range("G24").Value = 88
$ws1.Range("H24").Value = "https://show.bilibili.com/platform/detail.html?id=87527"
$ws1.Range("I24").Value = "//i1.hdslb.com/bfs/openplatform/202406/PgbS9Ic01718624347589.jpeg"
$ws1.Range("C25").Value = "杭州·夏之誓国乙only-日夜场"
$ws1.Range("D25").Value = "北干街道萧杭路689号 杭州时尚外滩艺术中心"
$ws1.Range("E25").Value = "2024.07.27 10:00-07.27 21:00"
$ws1.Range("F25").Value = 1655
$ws1.Range("G25").Value = 99
$ws1.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=83589"
$ws1.Range("I25").Value = "//i2.hdslb.com/bfs/openplatform/202405/99kWb2dy1714964533903.png"
$ws1.Range("F26").Value = 1655
$ws1.Range("F27").Value = 711
$ws1.Range("F28").Value = 909
$ws1.Range("F29").Value = 6
$ws1.Range("F30").Value = 1867
$ws1.Range("F32").Value = 2277
$ws1.Range("F34").Value = 108
$ws1.Range("F35").Value = 1452
$ws1.Range("F38").Value = 791
$ws1.Range("F39").Value = 407
$ws1.Range("F40").Value = 2962
$ws1.Range("F41").Value = 4083
$ws1.Range("F42").Value = 192
$ws1.Range("F47").Value = 859
$ws1.Range("F49").Value = 4080

# Sheet 3: 本地生活 (Local life)
$ws3.Range("F2").Value = 5192

# Sheet 4: 全部类型 (All types)
$ws4.Range("F4").Value = 857
$ws4.Range("B5").Value = "'2024-06-30"
$ws4.Range("C5").Value = "杭州·热血番ONLY1.0"
$ws4.Range("D5").Value = "康候圣街99号 顺丰创新中心"
$ws4.Range("E5").Value = "2024.06.30 10:00-06.30 17:00"
$ws4.Range("F5").Value = 857
$ws4.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=85042"
$ws4.Range("I5").Value = "//i2.hdslb.com/bfs/openplatform/202406/ckVVTuNj1717752114555.jpeg"
$ws4.Range("B6").Value = "'2024-07-04"
$ws4.Range("C6").Value = "杭州·乙游Camerata杭州2.0"
$ws4.Range("D6").Value = "杭海路601号江和美海洋广场1层 嘉宝丽酒店"
$ws4.Range("E6").Value = "2024.07.04 10:00-07.04 17:00"
$ws4.Range("F6").Value = 441
$ws4.Range("G6").Value = 68
$ws4.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=86777"
$ws4.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202406/aBDjuHlA1717403033570.jpeg"
$ws4.Range("C7").Value = "杭州·萌忧 动漫游戏嘉年华"
$ws4.Range("D7").Value = "康候圣街99号 顺丰创新中心"
$ws4.Range("E7").Value = "2024.07.06 10:30-07.06 17:00"
$ws4.Range("F7").Value = 889
$ws4.Range("G7").Value = 55
$ws4.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=87292"
$ws4.Range("I7").Value = "//i0.hdslb.com/bfs/openplatform/202406/s4DzfyI11718083528496.jpeg"
$ws4.Range("B8").Value = "'2024-07-06"
$ws4.Range("C8").Value = "杭州·黑执事only"
$ws4.Range("D8").Value = "大岭山路156号 爱丽芬城堡"
$ws4.Range("E8").Value = "2024.07.06 10:00-07.07 18:00"
$ws4.Range("F8").Value = 504
$ws4.Range("G8").Value = 160
$ws4.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=86414"
$ws4.Range("I8").Value = "//i1.hdslb.com/bfs/openplatform/202405/iP2cxk2w1716800288950.jpeg"
$ws4.Range("C9").Value = "杭州·【早鸟4折】《海上钢琴师1900》一生必听经典电影主题音乐会"
$ws4.Range("D9").Value = "湖墅南路138号 杭州浙话艺术剧院"
$ws4.Range("E9").Value = "2024.07.13 14:00-07.13 15:30"
$ws4.Range("F9").Value = 3
$ws4.Range("G9").Value = 72
$ws4.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=87485"
$ws4.Range("I9").Value = "//i0.hdslb.com/bfs/openplatform/202406/alBf4PMZ1718607067950.jpeg"
$ws4.Range("C10").Value = "杭州·【早鸟6折】《忱宴·渐渐被你吸引》热血动漫二次元ACG演唱会"
$ws4.Range("D10").Value = "湖墅南路136-138号 浙话艺术剧院"
$ws4.Range("E10").Value = "2024.07.13 19:30-07.13 21:30"
$ws4.Range("F10").Value = 21
$ws4.Range("G10").Value = 60
$ws4.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=85011"
$ws4.Range("I10").Value = "//i1.hdslb.com/bfs/openplatform/202404/2Gd8eLva1714379746993.jpeg"
$ws4.Range("C11").Value = "杭州·海上钢琴师—一生必听的电影名曲《泰坦尼克号》《花样年华》《海上钢琴师》"
$ws4.Range("D11").Value = "曙光路31号 浙江音乐厅"
$ws4.Range("E11").Value = "2024.07.13 19:30-07.13 21:00"
$ws4.Range("F11").Value = 8
$ws4.Range("G11").Value = 100
$ws4.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=85889"
$ws4.Range("I11").Value = "//i0.hdslb.com/bfs/openplatform/202405/52kxbBTh1716096935602.jpeg"
$ws4.Range("C12").Value = "杭州·草莓动漫节"
$ws4.Range("E12").Value = "2024.07.13 09:00-07.14 17:00"
$ws4.Range("F12").Value = 5496
$ws4.Range("G12").Value = 70
$ws4.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=84229"
$ws4.Range("I12").Value = "//i1.hdslb.com/bfs/openplatform/202406/czPRn1ve1718875288240.jpeg"
$ws4.Range("F13").Value = 7573
$ws4.Range("F15").Value = 1142
$ws4.Range("F16").Value = 893
$ws4.Range("F17").Value = 667
$ws4.Range("F18").Value = 220
$ws4.Range("F19").Value = 83
$ws4.Range("F23").Value = 1192
$ws4.Range("F24").Value = 110
$ws4.Range("F25").Value = 1655
$ws4.Range("F26").Value = 711
$ws4.Range("F27").Value = 909
$ws4.Range("F28").Value = 6
$ws4.Range("F29").Value = 1867
$ws4.Range("F31").Value = 2277
$ws4.Range("F39").Value = 407
$ws4.Range("F40").Value = 4083
$ws4.Range("F42").Value = 192
$ws4.Range("F47").Value = 859
$ws4.Range("F49").Value = 4080

Write-Host "Edit complete"
